$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 12 (start of group-2's rows), shifting
# the existing group-2 rows (12-17) down to 13-18.
$ws.Rows(12).Insert()

# The new row 12 becomes a blank placeholder row for "group-2" (only the
# groupID column is filled in, like the pre-existing blank row pattern
# used for group-2 before this edit).
$ws.Range("A12").Value = "group-2"

# The rows that used to belong to group-2 (now at 13-18) are re-labelled
# to a brand new group, "group-3".
$ws.Range("A13:A18").Value = "group-3"

$null = $ws.Range("N17").Select()
